# Auto-generated edit script applying the "add spin v2 xlsx" diff
# to xls/spin_the_movie.xlsx (rows 5-8 numeric columns + 4 label strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated test-run label strings (columns B and C), rows 5-8
$ws.Cells.Item(5, 2).Value = 'MT mem RD BW test, dura= 3.012559, GB/sec= 18.216385'
$ws.Cells.Item(5, 3).Value = 'MT mem RD BW test, dura= 3.012559, GB/sec= 18.216385'
$ws.Cells.Item(6, 2).Value = 'MT L3 RD BW test, dura= 3.000009, GB/sec= 65.738517'
$ws.Cells.Item(6, 3).Value = 'MT L3 RD BW test, dura= 3.000009, GB/sec= 65.738517'
$ws.Cells.Item(7, 2).Value = 'MT L2 RD BW test, dura= 3.000003, GB/sec= 133.187461'
$ws.Cells.Item(7, 3).Value = 'MT L2 RD BW test, dura= 3.000003, GB/sec= 133.187461'
$ws.Cells.Item(8, 2).Value = 'MT spin test, dura= 3.000000, Gops/sec= 3.857277'
$ws.Cells.Item(8, 3).Value = 'MT spin test, dura= 3.000000, Gops/sec= 3.857277'

# Updated numeric measurements, rows 5-8
# Row 5
$ws.Cells.Item(5, 4).Value = 828.367259
$ws.Cells.Item(5, 5).Value = 831.374318
$ws.Cells.Item(5, 7).Value = 400
$ws.Cells.Item(5, 8).Value = 49.8
$ws.Cells.Item(5, 9).Value = 2.4
$ws.Cells.Item(5, 10).Value = 2.4
$ws.Cells.Item(5, 11).Value = 2.4
$ws.Cells.Item(5, 12).Value = 2.4
$ws.Cells.Item(5, 13).Value = 18.2
$ws.Cells.Item(5, 14).Value = 9.800000000000001
$ws.Cells.Item(5, 16).Value = 13.9
$ws.Cells.Item(5, 18).Value = 18242.7
$ws.Cells.Item(5, 19).Value = 18327.7
$ws.Cells.Item(5, 20).Value = 85
$ws.Cells.Item(5, 21).Value = 72.3
$ws.Cells.Item(5, 22).Value = 0.6
$ws.Cells.Item(5, 23).Value = 0.6
$ws.Cells.Item(5, 24).Value = 64.09999999999999
$ws.Cells.Item(5, 25).Value = 63.9
$ws.Cells.Item(5, 26).Value = 6.6
$ws.Cells.Item(5, 27).Value = 6.6
$ws.Cells.Item(5, 28).Value = 20.8
$ws.Cells.Item(5, 29).Value = 20.8
$ws.Cells.Item(5, 32).Value = 11.2
$ws.Cells.Item(5, 33).Value = 11.1
$ws.Cells.Item(5, 34).Value = 3.8
$ws.Cells.Item(5, 35).Value = 3.8
$ws.Cells.Item(5, 36).Value = 0.6
$ws.Cells.Item(5, 37).Value = 0.8
$ws.Cells.Item(5, 38).Value = 3
$ws.Cells.Item(5, 39).Value = 2.3
$ws.Cells.Item(5, 40).Value = 0.6
$ws.Cells.Item(5, 41).Value = 0.6
$ws.Cells.Item(5, 43).Value = 0.03
$ws.Cells.Item(5, 44).Value = 0.35
$ws.Cells.Item(5, 45).Value = 0.35
$ws.Cells.Item(5, 46).Value = 0.03
$ws.Cells.Item(5, 47).Value = 0.03
$ws.Cells.Item(5, 48).Value = 88.2
$ws.Cells.Item(5, 49).Value = 87.7
$ws.Cells.Item(5, 50).Value = 88.40000000000001
$ws.Cells.Item(5, 51).Value = 88.09999999999999
$ws.Cells.Item(5, 52).Value = 0.7
$ws.Cells.Item(5, 53).Value = 1.1
$ws.Cells.Item(5, 54).Value = 0.42
$ws.Cells.Item(5, 55).Value = 0.43
$ws.Cells.Item(5, 56).Value = 0.07000000000000001
$ws.Cells.Item(5, 57).Value = 0.07000000000000001
$ws.Cells.Item(5, 58).Value = 0.07000000000000001
$ws.Cells.Item(5, 59).Value = 0.07000000000000001
$ws.Cells.Item(5, 60).Value = 0.04
$ws.Cells.Item(5, 61).Value = 0.04
$ws.Cells.Item(5, 62).Value = 0.04
$ws.Cells.Item(5, 63).Value = 0.04
$ws.Cells.Item(5, 64).Value = 0.01
$ws.Cells.Item(5, 66).Value = 0.07000000000000001
$ws.Cells.Item(5, 67).Value = 0.07000000000000001
$ws.Cells.Item(5, 68).Value = 0.12
$ws.Cells.Item(5, 69).Value = 0.13
$ws.Cells.Item(5, 72).Value = 0.16
$ws.Cells.Item(5, 73).Value = 0.16
$ws.Cells.Item(5, 74).Value = 0.17
$ws.Cells.Item(5, 75).Value = 0.17
$ws.Cells.Item(5, 76).Value = 0.09
$ws.Cells.Item(5, 77).Value = 0.09
$ws.Cells.Item(5, 78).Value = 0.09
$ws.Cells.Item(5, 79).Value = 0.09
$ws.Cells.Item(5, 81).Value = 0.02
$ws.Cells.Item(5, 82).Value = 0.17
$ws.Cells.Item(5, 83).Value = 0.17
$ws.Cells.Item(5, 84).Value = 0.3
$ws.Cells.Item(5, 85).Value = 0.3
$ws.Cells.Item(5, 86).Value = 0.01
# Row 6
$ws.Cells.Item(6, 4).Value = 831.380709
$ws.Cells.Item(6, 5).Value = 834.375218
$ws.Cells.Item(6, 8).Value = 56.8
$ws.Cells.Item(6, 13).Value = 22.6
$ws.Cells.Item(6, 18).Value = 519.1
$ws.Cells.Item(6, 20).Value = 5.8
$ws.Cells.Item(6, 21).Value = 0.2
$ws.Cells.Item(6, 24).Value = 29.3
$ws.Cells.Item(6, 25).Value = 29.2
$ws.Cells.Item(6, 28).Value = 43.4
$ws.Cells.Item(6, 29).Value = 43.1
$ws.Cells.Item(6, 32).Value = 25.8
$ws.Cells.Item(6, 33).Value = 25.5
$ws.Cells.Item(6, 36).Value = 0.1
$ws.Cells.Item(6, 39).Value = 0.7
$ws.Cells.Item(6, 43).Value = 0.03
$ws.Cells.Item(6, 44).Value = 1.18
$ws.Cells.Item(6, 48).Value = 66.8
$ws.Cells.Item(6, 49).Value = 66.7
$ws.Cells.Item(6, 51).Value = 66.09999999999999
$ws.Cells.Item(6, 53).Value = 1.3
$ws.Cells.Item(6, 55).Value = 1.23
$ws.Cells.Item(6, 67).Value = 0.19
$ws.Cells.Item(6, 79).Value = 0.28
$ws.Cells.Item(6, 80).Value = 0.02
$ws.Cells.Item(6, 81).Value = 0.01
$ws.Cells.Item(6, 82).Value = 0.5
$ws.Cells.Item(6, 83).Value = 0.5
# Row 7
$ws.Cells.Item(7, 4).Value = 834.388636
$ws.Cells.Item(7, 5).Value = 837.383139
$ws.Cells.Item(7, 7).Value = 399.7
$ws.Cells.Item(7, 8).Value = 61.5
$ws.Cells.Item(7, 13).Value = 23.1
$ws.Cells.Item(7, 16).Value = 16.7
$ws.Cells.Item(7, 18).Value = 493.5
$ws.Cells.Item(7, 19).Value = 498.6
$ws.Cells.Item(7, 20).Value = 5.1
$ws.Cells.Item(7, 24).Value = 0.2
$ws.Cells.Item(7, 25).Value = 0.1
$ws.Cells.Item(7, 26).Value = 1.5
$ws.Cells.Item(7, 27).Value = 1
$ws.Cells.Item(7, 28).Value = 41.2
$ws.Cells.Item(7, 29).Value = 40.4
$ws.Cells.Item(7, 32).Value = 17.6
$ws.Cells.Item(7, 33).Value = 17.2
$ws.Cells.Item(7, 34).Value = 25.8
$ws.Cells.Item(7, 35).Value = 26
$ws.Cells.Item(7, 38).Value = 0.1
$ws.Cells.Item(7, 40).Value = 0.2
$ws.Cells.Item(7, 44).Value = 2.37
$ws.Cells.Item(7, 45).Value = 2.38
$ws.Cells.Item(7, 47).Value = 0.01
$ws.Cells.Item(7, 48).Value = 32
$ws.Cells.Item(7, 49).Value = 31.6
$ws.Cells.Item(7, 50).Value = 34.2
$ws.Cells.Item(7, 51).Value = 34
$ws.Cells.Item(7, 52).Value = 1.7
$ws.Cells.Item(7, 53).Value = 2.1
$ws.Cells.Item(7, 54).Value = 2.47
$ws.Cells.Item(7, 55).Value = 2.5
$ws.Cells.Item(7, 56).Value = 0.38
$ws.Cells.Item(7, 57).Value = 0.38
$ws.Cells.Item(7, 58).Value = 0.39
$ws.Cells.Item(7, 59).Value = 0.39
$ws.Cells.Item(7, 60).Value = 0.21
$ws.Cells.Item(7, 61).Value = 0.21
$ws.Cells.Item(7, 62).Value = 0.21
$ws.Cells.Item(7, 63).Value = 0.21
$ws.Cells.Item(7, 66).Value = 0.39
$ws.Cells.Item(7, 68).Value = 0.64
$ws.Cells.Item(7, 69).Value = 0.64
$ws.Cells.Item(7, 72).Value = 0.98
$ws.Cells.Item(7, 73).Value = 0.99
$ws.Cells.Item(7, 74).Value = 1
$ws.Cells.Item(7, 75).Value = 1.01
$ws.Cells.Item(7, 76).Value = 0.55
$ws.Cells.Item(7, 77).Value = 0.55
$ws.Cells.Item(7, 78).Value = 0.55
$ws.Cells.Item(7, 79).Value = 0.55
$ws.Cells.Item(7, 82).Value = 1.01
$ws.Cells.Item(7, 83).Value = 1.02
$ws.Cells.Item(7, 84).Value = 1.65
$ws.Cells.Item(7, 85).Value = 1.66
$ws.Cells.Item(7, 87).Value = 0.02
# Row 8
$ws.Cells.Item(8, 4).Value = 837.392574
$ws.Cells.Item(8, 5).Value = 840.387074
$ws.Cells.Item(8, 8).Value = 65.5
$ws.Cells.Item(8, 13).Value = 23.1
$ws.Cells.Item(8, 14).Value = 13.1
$ws.Cells.Item(8, 18).Value = 482.8
$ws.Cells.Item(8, 19).Value = 487.3
$ws.Cells.Item(8, 20).Value = 4.5
$ws.Cells.Item(8, 28).Value = 0.1
$ws.Cells.Item(8, 29).Value = 0.1
$ws.Cells.Item(8, 32).Value = 0
$ws.Cells.Item(8, 34).Value = 0
$ws.Cells.Item(8, 35).Value = 0.1
$ws.Cells.Item(8, 39).Value = 0.1
$ws.Cells.Item(8, 42).Value = 2.87
$ws.Cells.Item(8, 43).Value = 2.84
$ws.Cells.Item(8, 44).Value = 0
$ws.Cells.Item(8, 45).Value = 0
$ws.Cells.Item(8, 48).Value = 11.4
$ws.Cells.Item(8, 49).Value = 11.4
$ws.Cells.Item(8, 50).Value = 8
$ws.Cells.Item(8, 51).Value = 8.199999999999999
$ws.Cells.Item(8, 52).Value = 0.7
$ws.Cells.Item(8, 53).Value = 0.8
$ws.Cells.Item(8, 54).Value = 3.31
$ws.Cells.Item(8, 55).Value = 3.31
$ws.Cells.Item(8, 56).Value = 0.5600000000000001
$ws.Cells.Item(8, 67).Value = 0.62
$ws.Cells.Item(8, 74).Value = 1.53
$ws.Cells.Item(8, 82).Value = 1.63
$ws.Cells.Item(8, 83).Value = 1.61
$ws.Cells.Item(8, 84).Value = 1.94
$ws.Cells.Item(8, 85).Value = 1.93
